$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.040.13"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "3.331.94"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").Value = "3.328.45"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.405"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "3.910.55"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "67.223.00"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "3.333.77"
$ws.Range("E18").Value = "  +1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "442.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "3.479.98"
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.512"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.195"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.67%  "
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("D41").Value = "2.826.29"
$ws.Range("E41").Value = "  +7.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.789"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0672"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "322.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.984"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
